{"js": "const body = context.document.body;\n\nasync function replaceOnce(oldText, newText) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst replacements = [\n  [\"2025-05-29 Thursday\", \"2025-05-30 Friday\"],\n  [\"747\u00d73=2241\", \"736\u00d74=2944\"],\n  [\"789\u00d79=7101\", \"676\u00d74=2704\"],\n  [\"782\u00d75=3910\", \"956\u00d75=4780\"],\n  [\"344\u00d72=688\", \"964\u00d75=4820\"],\n  [\"172\u00d72=344\", \"355\u00d79=3195\"],\n  [\"404\u00d76=2424\", \"723\u00d77=5061\"],\n  [\"191\u00d74=764\", \"616\u00d79=5544\"],\n  [\"406\u00d75=2030\", \"353\u00d77=2471\"],\n  [\"788\u00d75=3940\", \"644\u00d77=4508\"],\n  [\"275\u00d76=1650\", \"856\u00d76=5136\"],\n  [\"135\u00d72=270\", \"890\u00d72=1780\"],\n  [\"531\u00d77=3717\", \"497\u00d76=2982\"],\n  [\"782\u00d75=3910\", \"578\u00d76=3468\"],\n  [\"764\u00d79=6876\", \"780\u00d73=2340\"],\n  [\"572\u00d78=4576\", \"457\u00d75=2285\"],\n  [\"534\u00d79=4806\", \"191\u00d72=382\"],\n  [\"509\u00d75=2545\", \"276\u00d74=1104\"],\n  [\"756\u00d72=1512\", \"410\u00d76=2460\"],\n  [\"195\u00d73=585\", \"306\u00d79=2754\"],\n  [\"912\u00d77=6384\", \"748\u00d72=1496\"],\n  [\"550\u00d76=3300\", \"759\u00d77=5313\"],\n  [\"677\u00d72=1354\", \"551\u00d74=2204\"],\n  [\"543\u00d73=1629\", \"548\u00d73=1644\"],\n  [\"943\u00d74=3772\", \"233\u00d72=466\"],\n  [\"574\u00d74=2296\", \"963\u00d78=7704\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  await replaceOnce(oldText, newText);\n}\n", "ps1": "$d = $word.ActiveDocument\n$rng = $d.Content\n\nfunction Replace-NextOccurrence($range, $findText, $replaceText) {\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $ok = $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1)\n    if (-not $ok) {\n        throw \"Could not find text: $findText\"\n    }\n}\n\nReplace-NextOccurrence $rng \"2025-05-29 Thursday\" \"2025-05-30 Friday\"\nReplace-NextOccurrence $rng \"747\u00d73=2241\" \"736\u00d74=2944\"\nReplace-NextOccurrence $rng \"789\u00d79=7101\" \"676\u00d74=2704\"\nReplace-NextOccurrence $rng \"782\u00d75=3910\" \"956\u00d75=4780\"\nReplace-NextOccurrence $rng \"344\u00d72=688\" \"964\u00d75=4820\"\nReplace-NextOccurrence $rng \"172\u00d72=344\" \"355\u00d79=3195\"\nReplace-NextOccurrence $rng \"404\u00d76=2424\" \"723\u00d77=5061\"\nReplace-NextOccurrence $rng \"191\u00d74=764\" \"616\u00d79=5544\"\nReplace-NextOccurrence $rng \"406\u00d75=2030\" \"353\u00d77=2471\"\nReplace-NextOccurrence $rng \"788\u00d75=3940\" \"644\u00d77=4508\"\nReplace-NextOccurrence $rng \"275\u00d76=1650\" \"856\u00d76=5136\"\nReplace-NextOccurrence $rng \"135\u00d72=270\" \"890\u00d72=1780\"\nReplace-NextOccurrence $rng \"531\u00d77=3717\" \"497\u00d76=2982\"\nReplace-NextOccurrence $rng \"782\u00d75=3910\" \"578\u00d76=3468\"\nReplace-NextOccurrence $rng \"764\u00d79=6876\" \"780\u00d73=2340\"\nReplace-NextOccurrence $rng \"572\u00d78=4576\" \"457\u00d75=2285\"\nReplace-NextOccurrence $rng \"534\u00d79=4806\" \"191\u00d72=382\"\nReplace-NextOccurrence $rng \"509\u00d75=2545\" \"276\u00d74=1104\"\nReplace-NextOccurrence $rng \"756\u00d72=1512\" \"410\u00d76=2460\"\nReplace-NextOccurrence $rng \"195\u00d73=585\" \"306\u00d79=2754\"\nReplace-NextOccurrence $rng \"912\u00d77=6384\" \"748\u00d72=1496\"\nReplace-NextOccurrence $rng \"550\u00d76=3300\" \"759\u00d77=5313\"\nReplace-NextOccurrence $rng \"677\u00d72=1354\" \"551\u00d74=2204\"\nReplace-NextOccurrence $rng \"543\u00d73=1629\" \"548\u00d73=1644\"\nReplace-NextOccurrence $rng \"943\u00d74=3772\" \"233\u00d72=466\"\nReplace-NextOccurrence $rng \"574\u00d74=2296\" \"963\u00d78=7704\"\n"}
